# Add season record columns (Wins, Losses, Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row so the record columns line up with all data rows.
$lastRow = $ws.UsedRange.Rows.Count()

# --- Header row (row 1): new column headers, styled like the existing headers ---
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows: every team row gets the same season record ---
$ws.Range("AD2:AD" + $lastRow).Value = 75
$ws.Range("AE2:AE" + $lastRow).Value = 87
$ws.Range("AF2:AF" + $lastRow).Value = 0
